# Add a new row (TC row 17) to the "Test Steps" sheet, mirroring the
# existing "Login_03" test-case rows already present on "Test Cases".
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Test Steps")

$ws.Range("A17").Value = "Login_03"
$ws.Range("B17").Value = "TS_002"
$ws.Range("C17").Value = "Enter the password in Password field"
$ws.Range("D17").Value = "txt_Password"
$ws.Range("E17").Value = "input_Password"

# Leave the selection where the author left it after typing the new row.
$ws.Range("A14").Select()
